# Auto-generated edit script applying the cryptos.xlsx diff
# Updates Price (D) and Volume(1h) (E) columns, and for the two row-swaps
# also updates Coin (B) and Link (C) columns.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$colMap = @{ B = 2; C = 3; D = 4; E = 5 }

$data = @{
    2 = @{ D = "28.576.95"; E = "  +0.91%  " }
    3 = @{ D = "1.556.12"; E = "  -1.36%  " }
    4 = @{ D = "0.998"; E = "  -0.60%  " }
    5 = @{ E = "  -1.03%  " }
    6 = @{ E = "  -1.14%  " }
    7 = @{ D = "0.998"; E = "  -0.60%  " }
    8 = @{ D = "24.33"; E = "  +2.17%  " }
    9 = @{ D = "0.244"; E = "  -1.06%  " }
    10 = @{ E = "  -0.90%  " }
    11 = @{ E = "  -0.27%  " }
    12 = @{ D = "1.777.69"; E = "  -1.44%  " }
    13 = @{ D = "1.554.19"; E = "  -1.59%  " }
    14 = @{ D = "28.581.71"; E = "  +0.80%  " }
    15 = @{ D = "0.511"; E = "  -1.17%  " }
    16 = @{ E = "  -1.78%  " }
    17 = @{ D = "61.19"; E = "  -0.82%  " }
    18 = @{ D = "229.58"; E = "  -0.65%  " }
    19 = @{ D = "7.37"; E = "  -0.91%  " }
    21 = @{ E = "  -0.48%  " }
    22 = @{ E = "  -1.14%  " }
    23 = @{ E = "  -1.65%  " }
    24 = @{ D = "2.08"; E = "  +1.00%  " }
    25 = @{ D = "151.25"; E = "  -0.41%  " }
    26 = @{ E = "  -1.78%  " }
    27 = @{ E = "  -0.99%  " }
    28 = @{ E = "  -0.51%  " }
    29 = @{ E = "  -2.30%  " }
    30 = @{ E = "  -4.33%  " }
    31 = @{ E = "  -1.73%  " }
    32 = @{ D = "3.16"; E = "  -1.06%  " }
    33 = @{ D = "1.392.08"; E = "  +0.05%  " }
    34 = @{ E = "  -2.86%  " }
    35 = @{ D = "1.04"; E = "  -3.31%  " }
    36 = @{ E = "  -1.62%  " }
    37 = @{ B = "MXToken"; C = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"; D = "2.67"; E = "  +0.80%  " }
    38 = @{ B = "HuobiToken"; C = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"; D = "2.30"; E = "  -2.94%  " }
    40 = @{ B = "ImmutableX"; C = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"; D = "0.518"; E = "  -0.31%  " }
    41 = @{ B = "RenderToken"; C = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"; D = "1.94"; E = "  +2.38%  " }
    42 = @{ E = "  -0.49%  " }
    43 = @{ D = "0.776"; E = "  -1.23%  " }
    44 = @{ D = "0.0462"; E = "  +1.27%  " }
    45 = @{ D = "64.24"; E = "  +2.82%  " }
    46 = @{ E = "  -2.09%  " }
    47 = @{ D = "1.690.39"; E = "  -1.44%  " }
    48 = @{ D = "0.866"; E = "  -6.46%  " }
    49 = @{ D = "43.72"; E = "  +6.36%  " }
    50 = @{ D = "85.24"; E = "  -0.37%  " }
    51 = @{ D = "0.0₆0101"; E = "  -0.82%  " }
}

foreach ($row in $data.Keys) {
    $cols = $data[$row]
    foreach ($col in $cols.Keys) {
        $colIdx = $colMap[$col]
        $val = $cols[$col]
        $cell = $ws.Cells.Item($row, $colIdx)
        if ($col -eq "D") {
            # Force text storage so purely-numeric-looking price strings
            # (e.g. "0.998") are not reinterpreted as numbers by Excel.
            $cell.NumberFormat = "@"
            $cell.Value = $val
            $cell.ClearFormats()
        } else {
            $cell.Value = $val
        }
    }
}

Write-Output "Applied $($data.Count) row updates"
